# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count) figures and one event title
# across the "展览", "演出", "本地生活" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 269
$ws.Range("F4").Value = 1092
$ws.Range("F5").Value = 2573
$ws.Range("F6").Value = 227
$ws.Range("F7").Value = 679
$ws.Range("F8").Value = 56
$ws.Range("F9").Value = 237
$ws.Range("C11").Value = "广州·wio jump 同人only4.0"
$ws.Range("F11").Value = 686
$ws.Range("F12").Value = 84
$ws.Range("F13").Value = 109
$ws.Range("F14").Value = 1482
$ws.Range("F15").Value = 241
$ws.Range("F16").Value = 45

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 17
$ws.Range("F19").Value = 50

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6346
$ws.Range("F4").Value = 2010
$ws.Range("F5").Value = 243

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types) -- combined listing, same data repeated
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6346
$ws.Range("F4").Value = 2010
$ws.Range("F5").Value = 243
$ws.Range("F11").Value = 269
$ws.Range("F12").Value = 1092
$ws.Range("F16").Value = 2573
$ws.Range("F18").Value = 227
$ws.Range("F19").Value = 17
$ws.Range("F22").Value = 679
$ws.Range("F23").Value = 56
$ws.Range("F24").Value = 237
$ws.Range("C27").Value = "广州·wio jump 同人only4.0"
$ws.Range("F27").Value = 686
$ws.Range("F28").Value = 84
$ws.Range("F29").Value = 109
$ws.Range("F31").Value = 1482
$ws.Range("F32").Value = 242
$ws.Range("F35").Value = 45
$ws.Range("F39").Value = 50
